$wb = $excel.ActiveWorkbook

$handbackMsg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/e4bb53bbbe500451e4bbec813100bbcda7620194/e2e/9938cbe1-33c3-4a1f-831a-66db315a3c73.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/c5b40d8aa0a7ce3a6cb95a69f5a905a60e4d5aae/e2e/9938cbe1-33c3-4a1f-831a-66db315a3c73.md."

# ---- zh-cn sheet: row 7 now has a completed handback ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I7").Value = "9938cbe1-33c3-4a1f-831a-66db315a3c73.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/6b1f6e7cf389ca5c5ba27f1d0c5b2d1f6d76dba2/e2e/9938cbe1-33c3-4a1f-831a-66db315a3c73.md", [System.Type]::Missing, [System.Type]::Missing, "9938cbe1-33c3-4a1f-831a-66db315a3c73.md")
$wsZh.Range("I7").Font.Underline = 2
$wsZh.Range("I7").Font.Color = 15570276

$wsZh.Range("J7").Value = "9938cbe1-33c3-4a1f-831a-66db315a3c73.b3c4458a561d8fc5f8ee25ec12c2705e07a0cbb8.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-13 07:09:16"
$wsZh.Range("P7").Value = $handbackMsg

# ---- de-de sheet: row 7 now has a completed handback ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I7").Value = "9938cbe1-33c3-4a1f-831a-66db315a3c73.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6b1f6e7cf389ca5c5ba27f1d0c5b2d1f6d76dba2/e2e/9938cbe1-33c3-4a1f-831a-66db315a3c73.md", [System.Type]::Missing, [System.Type]::Missing, "9938cbe1-33c3-4a1f-831a-66db315a3c73.md")
$wsDe.Range("I7").Font.Underline = 2
$wsDe.Range("I7").Font.Color = 15570276

$wsDe.Range("J7").Value = "9938cbe1-33c3-4a1f-831a-66db315a3c73.b3c4458a561d8fc5f8ee25ec12c2705e07a0cbb8.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-13 07:09:26"
$wsDe.Range("P7").Value = $handbackMsg
